$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 and E2 hold numeric-looking identifiers that must stay text (as in the
# original inlineStr cells), so force a text number format before writing
# them to avoid Excel auto-converting them to numbers.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "400011185825004"
$ws.Range("C2").Value = "ΤΡΑΚΑΔΑΣ Α.Ε."
$ws.Range("D2").Value = "8Μ0ΤΔΑ"
$ws.Range("E2").Value = "8961"
$ws.Range("I2").Value = "22,61"
$ws.Range("J2").Value = "5,43"
$ws.Range("K2").Value = "28,04"
